# The workbook originally has a single sheet ("Sheet1") holding the OD600
# bacterial-growth data. The edit renames that sheet to "Sheet2" and
# inserts a brand new "Sheet1" in front of it that just explains that the
# real data now lives on the next sheet.

$wb = $excel.ActiveWorkbook

# Rename the existing (data) worksheet to "Sheet2".
$dataSheet = $wb.Worksheets.Item(1)
$dataSheet.Name = "Sheet2"

# Insert a new worksheet before it; Excel names it "Sheet1" automatically.
$introSheet = $wb.Worksheets.Add($wb.Sheets("Sheet2"))
$introSheet.Name = "Sheet1"

# Fill in the little explanatory note on the new intro sheet.
$wb.Sheets("Sheet1").Range("A1").Value = "Здесь могло бы быть описание данных."
$wb.Sheets("Sheet1").Range("A2").Value = "А сами данные в этом файле на следующем листе ;)"
$wb.Sheets("Sheet1").Range("A3").Select()

# Restore the data sheet's own scroll position/selection (it used to be the
# active sheet with F18 selected; now it just keeps its own view state).
$wb.Sheets("Sheet2").Activate()
$excel.ActiveWindow.ScrollRow = 35
$excel.ActiveWindow.ScrollColumn = 1
$wb.Sheets("Sheet2").Range("F10").Select()

# The new intro sheet is the one the user actually ends up looking at.
$wb.Sheets("Sheet1").Activate()
